$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so that numeric-looking
# strings (e.g. "506.53") are preserved exactly as text instead of being
# converted into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.183.22"
$ws.Range("D3").Value = "2.586.92"
$ws.Range("E3").Value = "  +5.93%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "506.53"
$ws.Range("E5").Value = "  +3.00%  "
$ws.Range("D6").Value = "155.29"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.581"
$ws.Range("E8").Value = "  -5.45%  "
$ws.Range("D9").Value = "2.625.81"
$ws.Range("E9").Value = "  +6.27%  "
$ws.Range("D10").Value = "6.47"
$ws.Range("E10").Value = "  +3.05%  "
$ws.Range("E11").Value = "  +2.53%  "
$ws.Range("D12").Value = "0.342"
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").Value = "3.064.82"
$ws.Range("E14").Value = "  +6.63%  "
$ws.Range("D15").Value = "60.306.39"
$ws.Range("E15").Value = "  +4.77%  "
$ws.Range("D16").Value = "21.67"
$ws.Range("E16").Value = "  +4.13%  "
$ws.Range("E17").Value = "  +3.58%  "
$ws.Range("D18").Value = "2.622.17"
$ws.Range("E18").Value = "  +6.45%  "
$ws.Range("D19").Value = "4.79"
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("D20").Value = "342.19"
$ws.Range("E20").Value = "  +4.84%  "
$ws.Range("D21").Value = "10.40"
$ws.Range("E21").Value = "  +3.26%  "
$ws.Range("D22").Value = "6.09"
$ws.Range("E22").Value = "  +2.32%  "
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "60.03"
$ws.Range("E24").Value = "  +2.83%  "
$ws.Range("E25").Value = "  +4.31%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "0.166"
$ws.Range("E26").Value = "  +3.05%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.708.34"
$ws.Range("E27").Value = "  +5.74%  "
$ws.Range("D28").Value = "0.989"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").Value = "0.0₃0859"
$ws.Range("E29").Value = "  +6.55%  "
$ws.Range("D30").Value = "7.51"
$ws.Range("E30").Value = "  +2.10%  "
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "156.07"
$ws.Range("E32").Value = "  +3.47%  "
$ws.Range("D33").Value = "19.31"
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("D35").Value = "5.73"
$ws.Range("E35").Value = "  +6.99%  "
$ws.Range("D36").Value = "3.99"
$ws.Range("E36").Value = "  +4.88%  "
$ws.Range("E37").Value = "  +5.15%  "
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").Value = "0.859"
$ws.Range("E38").Value = "  +28.80%  "
$ws.Range("D39").Value = "3.79"
$ws.Range("E39").Value = "  +6.52%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.48"
$ws.Range("E40").Value = "  +5.92%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "0.845"
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("D42").Value = "300.65"
$ws.Range("E42").Value = "  +6.54%  "
$ws.Range("D43").Value = "35.59"
$ws.Range("E43").Value = "  +3.76%  "
$ws.Range("E44").Value = "  +3.36%  "
$ws.Range("D45").Value = "0.0568"
$ws.Range("E45").Value = "  +5.77%  "
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").Value = "0.993"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "19.84"
$ws.Range("E48").Value = "  +10.45%  "
$ws.Range("D49").Value = "4.97"
$ws.Range("E49").Value = "  +5.81%  "
$ws.Range("D50").Value = "2.051.86"
$ws.Range("E50").Value = "  +8.29%  "
$ws.Range("D51").Value = "0.0233"
$ws.Range("E51").Value = "  +1.26%  "
